# Update "想去人数" (want-to-go count) values in the F column
# on both the "展览" and "全部类型" sheets to match the newly
# scraped data snapshot.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 1479
    $ws.Range("F3").Value = 3106
    $ws.Range("F4").Value = 45
    $ws.Range("F5").Value = 756
}
